$d = $word.ActiveDocument

# 1) Heading1 'ContextFreeSQL' -- wrap existing run with proofErr spellStart/spellEnd
$rng = $d.Content
$rng.Find.Execute("ContextFreeSQL") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ContextFreeSQL</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p.Range.InsertXML($xml)

# 2) Big restructuring of the 'schema compare' bullet list:
#    replace paragraphs from 'the code generation.' through 'The full load from connstr given in command line'
#    with the new RN/core/back-to-schema-loader/cleanup/views/full-load paragraphs
$rng = $d.Content
$rng.Find.Execute("the code generation.") | Out-Null
$pStart = $rng.Paragraphs(1)
$rng2 = $d.Content
$rng2.Find.Execute("The full load from connstr given in command line") | Out-Null
$pEnd = $rng2.Paragraphs(1)
$start = $pStart.Range.Start
$end = $pEnd.Range.End
$r = $d.Range($start, $end)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>RN</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>generation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> Work</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> on those variables (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>var_prefix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">). </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>where</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> are they getting values from in the .net code, move that into this. Then, dump it all into a file at the end, for now just hard code a file name</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>The core</w:t></w:r><w:r><w:t xml:space="preserve">! </w:t></w:r><w:r><w:t xml:space="preserve">go proc by proc on .net, just </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mimick</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> that</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>Back to schema loader:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Cleanup: </w:t></w:r><w:r><w:t xml:space="preserve">What about </w:t></w:r><w:r><w:t xml:space="preserve">column </w:t></w:r><w:r><w:t>defaults? Are they not in an already existing query?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Views, stored </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>procs</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The full load from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>connstr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> given in command line</w:t></w:r></w:p>'
$r.InsertXML($xml)

# 3) 'maintain all 3 options...' -- proofErr run split (text unchanged)
$rng = $d.Content
$rng.Find.Execute("maintain all 3 options") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">maintain all 3 </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>options, or</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> maybe move down to 2. Clearer docs in command line</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# 4) 'Full catching up with old product' heading -- proofErr run split (text unchanged)
$rng = $d.Content
$rng.Find.Execute("Full catching up with old product") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t xml:space="preserve">Full catching up with old </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>product</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$p.Range.InsertXML($xml)

# 5) 'Any other feature there?...' -- proofErr run split (text unchanged)
$rng = $d.Content
$rng.Find.Execute("Any other feature there") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Any other </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>feature</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> there? See again all the options on the GUI</w:t></w:r><w:r><w:t>, schema and data</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# 6) 'Setting data while ignoring IDENTITY...' -- proofErr run split (text unchanged)
$rng = $d.Content
$rng.Find.Execute("Setting data while ignoring IDENTITY") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Setting data while ignoring IDENTITY. Can either enforce numbers on IDENTITY fields (PG </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>most</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> allow it somewhere) or decide that they’re not important</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# 7) 'Smallie: ability to give...' -- proofErr run split (text unchanged)
$rng = $d.Content
$rng.Find.Execute("Smallie: ability to give") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Smallie: ability to give an AWS secret name in command </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>line ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> instead of full </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>connstr</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$p.Range.InsertXML($xml)

# 8) 'Loading from a SQL script, (sqlglot)...' -- proofErr run split (text unchanged)
$rng = $d.Content
$rng.Find.Execute("Loading from a SQL script") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Loading from a SQL script, </w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sqlglot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t xml:space="preserve">there can be CREATE, ALTER and DML data in there (and of course the scripting option to not remove </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>whats</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> already there, in case we’re loading from a partial script)</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# 9) 'Data script from pg backup?' -- proofErr run split (text unchanged)
$rng = $d.Content
$rng.Find.Execute("Data script from pg backup") | Out-Null
$p = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Data script from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> backup?</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

Write-Output "edit complete"
